# Apply odds updates for Jogos_da_Semana_FlashScore_2025-05-15.xlsx
# (commit: "Atualizando o arquivo XLSX")
# Updates numeric odds cells across several rows of Sheet1 to match the
# refreshed FlashScore snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = 2.03
$ws.Range("O3").Value = 1.87
$ws.Range("G7").Value = 1.6
$ws.Range("H7").Value = 4.1
$ws.Range("I7").Value = 5.25
$ws.Range("R7").Value = 1.8
$ws.Range("S7").Value = 1.91
$ws.Range("T7").Value = 7.5
$ws.Range("U7").Value = 8
$ws.Range("W7").Value = 12
$ws.Range("AA7").Value = 8
$ws.Range("AB7").Value = 17
$ws.Range("AG7").Value = 17
$ws.Range("AI7").Value = 41
$ws.Range("J13").Value = 1.1
$ws.Range("K13").Value = 7
$ws.Range("L13").Value = 1.5
$ws.Range("M13").Value = 2.5
$ws.Range("N13").Value = 2.5
$ws.Range("O13").Value = 1.5
$ws.Range("P13").Value = 1.53
$ws.Range("Q13").Value = 2.38
$ws.Range("R13").Value = 2.1
$ws.Range("S13").Value = 1.67
$ws.Range("Y13").Value = 41
$ws.Range("AE13").Value = 8
$ws.Range("G15").Value = 3.6
$ws.Range("H15").Value = 3.15
$ws.Range("I15").Value = 1.95
$ws.Range("N15").Value = 2.02
$ws.Range("O15").Value = 1.62
$ws.Range("P15").Value = 1.39
$ws.Range("Q15").Value = 2.45
$ws.Range("T15").Value = 7.9
$ws.Range("U15").Value = 15
$ws.Range("V15").Value = 10.5
$ws.Range("W15").Value = 40
$ws.Range("X15").Value = 28
$ws.Range("Y15").Value = 35
$ws.Range("Z15").Value = 8
$ws.Range("AA15").Value = 5.4
$ws.Range("AB15").Value = 12.5
$ws.Range("AC15").Value = 60
$ws.Range("AD15").Value = 450
$ws.Range("AE15").Value = 5.7
$ws.Range("AF15").Value = 7.5
$ws.Range("AG15").Value = 7.2
$ws.Range("AH15").Value = 14
$ws.Range("AI15").Value = 13.5
$ws.Range("H16").Value = 3.45
$ws.Range("O16").Value = 1.87
$ws.Range("U16").Value = 10
$ws.Range("Z16").Value = 11.25
$ws.Range("AB16").Value = 11
$ws.Range("AE16").Value = 8.25
$ws.Range("AJ16").Value = 22
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 9
$ws.Range("N17").Value = 2.2
$ws.Range("O17").Value = 1.65
$ws.Range("N18").Value = 2.3
$ws.Range("O18").Value = 1.6
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 3
$ws.Range("I20").Value = 2.35
$ws.Range("L20").Value = 1.47
$ws.Range("M20").Value = 2.32
$ws.Range("N20").Value = 2.35
$ws.Range("O20").Value = 1.47
$ws.Range("P20").Value = 1.55
$ws.Range("Q20").Value = 2.15
$ws.Range("R20").Value = 2.02
$ws.Range("S20").Value = 1.62
$ws.Range("T20").Value = 7.1
$ws.Range("U20").Value = 13.5
$ws.Range("V20").Value = 11.75
$ws.Range("W20").Value = 37
$ws.Range("X20").Value = 32
$ws.Range("Y20").Value = 50
$ws.Range("Z20").Value = 6.7
$ws.Range("AA20").Value = 6
$ws.Range("AB20").Value = 18.5
$ws.Range("AC20").Value = 120
$ws.Range("AE20").Value = 6.1
$ws.Range("AF20").Value = 10
$ws.Range("AG20").Value = 9.75
$ws.Range("AH20").Value = 24
$ws.Range("AI20").Value = 23
$ws.Range("AJ20").Value = 45
$ws.Range("G21").Value = 2.3
$ws.Range("I21").Value = 3.15
$ws.Range("L21").Value = 1.47
$ws.Range("M21").Value = 2.35
$ws.Range("N21").Value = 2.32
$ws.Range("O21").Value = 1.47
$ws.Range("P21").Value = 1.55
$ws.Range("Q21").Value = 2.15
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 1.65
$ws.Range("T21").Value = 6.1
$ws.Range("U21").Value = 10
$ws.Range("V21").Value = 9.5
$ws.Range("W21").Value = 23
$ws.Range("X21").Value = 22
$ws.Range("Y21").Value = 40
$ws.Range("Z21").Value = 6.7
$ws.Range("AB21").Value = 18
$ws.Range("AC21").Value = 110
$ws.Range("AD21").Value = 900
$ws.Range("AE21").Value = 7.3
$ws.Range("AF21").Value = 14.5
$ws.Range("AG21").Value = 12
$ws.Range("AI21").Value = 35
$ws.Range("AJ21").Value = 55
$ws.Range("G24").Value = 3.5
$ws.Range("J26").Value = 1.03
$ws.Range("K26").Value = 10.5
$ws.Range("L26").Value = 1.14
$ws.Range("M26").Value = 5
$ws.Range("N26").Value = 1.5
$ws.Range("O26").Value = 2.5
$ws.Range("G27").Value = 11
$ws.Range("I27").Value = 1.17
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = 1.02
$ws.Range("T27").Value = 34
$ws.Range("V27").Value = 29
$ws.Range("W27").Value = 126
$ws.Range("AG27").Value = 9.5
$ws.Range("AJ27").Value = 21
$ws.Range("G30").Value = 1.95
$ws.Range("H30").Value = 3.75
$ws.Range("I30").Value = 3.5
$ws.Range("R30").Value = 1.57
$ws.Range("S30").Value = 2.25
$ws.Range("U30").Value = 11
$ws.Range("X30").Value = 15
$ws.Range("AB30").Value = 12
$ws.Range("AD30").Value = 126
$ws.Range("AE30").Value = 13
$ws.Range("AG30").Value = 12
$ws.Range("AI30").Value = 26
$ws.Range("G32").Value = 1.75
$ws.Range("H32").Value = 3.8
$ws.Range("I32").Value = 4.2
$ws.Range("T32").Value = 8
$ws.Range("U32").Value = 9
$ws.Range("AB32").Value = 15
$ws.Range("AD32").Value = 201
$ws.Range("AF32").Value = 23
$ws.Range("AG32").Value = 15
$ws.Range("AI32").Value = 34
$ws.Range("G33").Value = 6.25
$ws.Range("I33").Value = 1.5
$ws.Range("R33").Value = 2
$ws.Range("S33").Value = 1.75
$ws.Range("X33").Value = 51
$ws.Range("AC33").Value = 67
$ws.Range("AD33").Value = 401
$ws.Range("AH33").Value = 10
$ws.Range("G36").Value = 2.35
$ws.Range("I36").Value = 3.2
$ws.Range("J36").Value = 1.08
$ws.Range("K36").Value = 8
$ws.Range("L36").Value = 1.36
$ws.Range("M36").Value = 3
$ws.Range("N36").Value = 2.15
$ws.Range("O36").Value = 1.67
$ws.Range("W36").Value = 21
$ws.Range("AD36").Value = 301
$ws.Range("AI36").Value = 29
$ws.Range("AJ36").Value = 41
